# Case_1_86 res_bus vm_pu.xlsx update: slack bus vm_pu set from 1.05 to 1.02 pu (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.018363534966368
$ws.Range("D2").Value = 1.028858935399628
$ws.Range("E2").Value = 1.019178495281031
$ws.Range("F2").Value = 1.028562777917307
$ws.Range("I2").Value = 1.030857377885566
$ws.Range("J2").Value = 1.023572392774784
$ws.Range("K2").Value = 1.031674528738413
$ws.Range("L2").Value = 1.02202251283487
$ws.Range("M2").Value = 1.031379232221555
$ws.Range("N2").Value = 1.025025982977921
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022436101562407
$ws.Range("D3").Value = 1.029453748545057
$ws.Range("E3").Value = 1.022759059100701
$ws.Range("F3").Value = 1.032139263568543
$ws.Range("I3").Value = 1.031145081514337
$ws.Range("J3").Value = 1.027264935404385
$ws.Range("K3").Value = 1.032077940646395
$ws.Range("L3").Value = 1.025401434221294
$ws.Range("M3").Value = 1.034756233422951
$ws.Range("N3").Value = 1.028723769441594
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025040311581058
$ws.Range("D4").Value = 1.029834334970225
$ws.Range("E4").Value = 1.025047614207758
$ws.Range("F4").Value = 1.034422139585006
$ws.Range("I4").Value = 1.03132545639031
$ws.Range("J4").Value = 1.029624304464952
$ws.Range("K4").Value = 1.0323339926566
$ws.Range("L4").Value = 1.027559595598305
$ws.Range("M4").Value = 1.036910101089983
$ws.Range("N4").Value = 1.03108648907685
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026127953197794
$ws.Range("D5").Value = 1.029993319767554
$ws.Range("E5").Value = 1.02600317184532
$ws.Range("F5").Value = 1.035374573342712
$ws.Range("I5").Value = 1.031399917506426
$ws.Range("J5").Value = 1.030609248237686
$ws.Range("K5").Value = 1.032440459226045
$ws.Range("L5").Value = 1.028460347792354
$ws.Range("M5").Value = 1.037808309127689
$ws.Range("N5").Value = 1.032072831582716
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026310159945959
$ws.Range("D6").Value = 1.030019954994815
$ws.Range("E6").Value = 1.026163236462536
$ws.Range("F6").Value = 1.035534070019931
$ws.Range("I6").Value = 1.031412340097553
$ws.Range("J6").Value = 1.030774224664462
$ws.Range("K6").Value = 1.032458266798297
$ws.Range("L6").Value = 1.028611210824967
$ws.Range("M6").Value = 1.037958701398689
$ws.Range("N6").Value = 1.032238042294942
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02505487253426
$ws.Range("D7").Value = 1.029836463296352
$ws.Range("E7").Value = 1.025060407861569
$ws.Range("F7").Value = 1.034434894401833
$ws.Range("I7").Value = 1.031326456697921
$ws.Range("J7").Value = 1.029637492274459
$ws.Range("K7").Value = 1.032335419875838
$ws.Range("L7").Value = 1.027571656901599
$ws.Range("M7").Value = 1.036922131313245
$ws.Range("N7").Value = 1.031099695614559
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.019746456540988
$ws.Range("D8").Value = 1.029060854403031
$ws.Range("E8").Value = 1.020394564003887
$ws.Range("F8").Value = 1.02977809268843
$ws.Range("I8").Value = 1.030955817905421
$ws.Range("J8").Value = 1.024826651452261
$ws.Range("K8").Value = 1.031811904650622
$ws.Range("L8").Value = 1.02317041380659
$ws.Range("M8").Value = 1.032527110780854
$ws.Range("N8").Value = 1.026282022846583
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.010142865609482
$ws.Range("D9").Value = 1.027660506246793
$ws.Range("E9").Value = 1.01194537372912
$ws.Range("F9").Value = 1.021321945548974
$ws.Range("I9").Value = 1.030257556070363
$ws.Range("J9").Value = 1.016108923881369
$ws.Range("K9").Value = 1.030850530692611
$ws.Range("L9").Value = 1.015188470630364
$ws.Range("M9").Value = 1.02453316042736
$ws.Range("N9").Value = 1.017551915102603
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.003555306664078
$ws.Range("D10").Value = 1.026703279827848
$ws.Range("E10").Value = 1.006144317794328
$ws.Range("F10").Value = 1.015501429134079
$ws.Range("I10").Value = 1.029760494865131
$ws.Range("J10").Value = 1.010119429586317
$ws.Range("K10").Value = 1.030182426680671
$ws.Range("L10").Value = 1.009700092204516
$ws.Range("M10").Value = 1.019021901790512
$ws.Range("N10").Value = 1.011553915038647
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.000654802644803
$ws.Range("D11").Value = 1.026282926516313
$ws.Range("E11").Value = 1.003588867804939
$ws.Range("F11").Value = 1.012934168737142
$ws.Range("I11").Value = 1.029537493196719
$ws.Range("J11").Value = 1.007479997200861
$ws.Range("K11").Value = 1.029886430674151
$ws.Range("L11").Value = 1.007280435599082
$ws.Range("M11").Value = 1.016588941637062
$ws.Range("N11").Value = 1.008910734356457
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 0.9995698392296042
$ws.Range("D12").Value = 1.026125884061977
$ws.Range("E12").Value = 1.002632791080093
$ws.Range("F12").Value = 1.011973208425337
$ws.Range("I12").Value = 1.029453467136003
$ws.Range("J12").Value = 1.006492353122088
$ws.Range("K12").Value = 1.029775454762114
$ws.Range("L12").Value = 1.006374869786914
$ws.Range("M12").Value = 1.015677935751724
$ws.Range("N12").Value = 1.007921687709808
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 0.9998029170191616
$ws.Range("D13").Value = 1.026159611515278
$ws.Range("E13").Value = 1.002838189015253
$ws.Range("F13").Value = 1.012179676187265
$ws.Range("I13").Value = 1.029471545441219
$ws.Range("J13").Value = 1.006704539451712
$ws.Range("K13").Value = 1.029799306453944
$ws.Range("L13").Value = 1.006569429558277
$ws.Range("M13").Value = 1.015873684802873
$ws.Range("N13").Value = 1.008134175368359
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.000565276200375
$ws.Range("D14").Value = 1.026269963927012
$ws.Range("E14").Value = 1.0035099802644
$ws.Range("F14").Value = 1.012854887579912
$ws.Range("I14").Value = 1.029530572073723
$ws.Range("J14").Value = 1.007398507975828
$ws.Range("K14").Value = 1.029877278515026
$ws.Range("L14").Value = 1.007205721789445
$ws.Range("M14").Value = 1.016513788249618
$ws.Range("N14").Value = 1.008829129407377
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.001033974158023
$ws.Range("D15").Value = 1.02633783512522
$ws.Range("E15").Value = 1.003922972837187
$ws.Range("F15").Value = 1.013269921962362
$ws.Range("I15").Value = 1.02956678140897
$ws.Range("J15").Value = 1.007825114905176
$ws.Range("K15").Value = 1.029925182510508
$ws.Range("L15").Value = 1.007596852013053
$ws.Range("M15").Value = 1.016907200987294
$ws.Range("N15").Value = 1.009256342167497
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.003746760223723
$ws.Range("D16").Value = 1.026731051734019
$ws.Range("E16").Value = 1.006312969557183
$ws.Range("F16").Value = 1.015670794406423
$ws.Range("I16").Value = 1.029775128854896
$ws.Range("J16").Value = 1.010293603517454
$ws.Range("K16").Value = 1.030201927825205
$ws.Range("L16").Value = 1.009859741187768
$ws.Range("M16").Value = 1.019182363117069
$ws.Range("N16").Value = 1.011728336316745
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.005435301172917
$ws.Range("D17").Value = 1.026976119320999
$ws.Range("E17").Value = 1.007800266019449
$ws.Range("F17").Value = 1.017164015043088
$ws.Range("I17").Value = 1.029903721415655
$ws.Range("J17").Value = 1.011829487139533
$ws.Range("K17").Value = 1.030373712656494
$ws.Range("L17").Value = 1.011267421533571
$ws.Range("M17").Value = 1.02059684077031
$ws.Range("N17").Value = 1.013266401069735
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00641558848173
$ws.Range("D18").Value = 1.027118498541419
$ws.Range("E18").Value = 1.008663601303077
$ws.Range("F18").Value = 1.018030478480373
$ws.Range("I18").Value = 1.029977979319057
$ws.Range("J18").Value = 1.012720933122862
$ws.Range("K18").Value = 1.030473266750314
$ws.Range("L18").Value = 1.012084356109977
$ws.Range("M18").Value = 1.021417411591735
$ws.Range("N18").Value = 1.014159113008611
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.006749070062402
$ws.Range("D19").Value = 1.027166951092453
$ws.Range("E19").Value = 1.008957277005795
$ws.Range("F19").Value = 1.018325164937859
$ws.Range("I19").Value = 1.030003173239246
$ws.Range("J19").Value = 1.013024155352339
$ws.Range("K19").Value = 1.030507103377382
$ws.Range("L19").Value = 1.012362216500718
$ws.Range("M19").Value = 1.021696455780973
$ws.Range("N19").Value = 1.014462765848434
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.005254616153952
$ws.Range("D20").Value = 1.026949884444738
$ws.Range("E20").Value = 1.0076411277975
$ws.Range("F20").Value = 1.017004274986595
$ws.Range("I20").Value = 1.029890002198405
$ws.Range("J20").Value = 1.011665159797827
$ws.Range("K20").Value = 1.030355348655876
$ws.Range("L20").Value = 1.011116821357669
$ws.Range("M20").Value = 1.020445545204286
$ws.Range("N20").Value = 1.013101840364355
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.000340992916277
$ws.Range("D21").Value = 1.026237493051955
$ws.Range("E21").Value = 1.003312346792953
$ws.Range("F21").Value = 1.012656260574973
$ws.Range("I21").Value = 1.029513223357308
$ws.Range("J21").Value = 1.007194354249149
$ws.Range("K21").Value = 1.029854346304143
$ws.Range("L21").Value = 1.007018539838946
$ws.Range("M21").Value = 1.016325497656357
$ws.Range("N21").Value = 1.00862468575899
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9972075102887076
$ws.Range("D22").Value = 1.025784338744477
$ws.Range("E22").Value = 1.00055075638638
$ws.Range("F22").Value = 1.009879710007049
$ws.Range("I22").Value = 1.029269409598919
$ws.Range("J22").Value = 1.004341304478934
$ws.Range("K22").Value = 1.029533375806283
$ws.Range("L22").Value = 1.004402291311703
$ws.Range("M22").Value = 1.013692685785617
$ws.Range("N22").Value = 1.005767584330851
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.99887293614132
$ws.Range("D23").Value = 1.026025069629904
$ws.Range("E23").Value = 1.002018624015618
$ws.Range("F23").Value = 1.011355776018105
$ws.Range("I23").Value = 1.029399324616227
$ws.Range("J23").Value = 1.005857866295638
$ws.Range("K23").Value = 1.029704102225321
$ws.Range("L23").Value = 1.005793067118073
$ws.Range("M23").Value = 1.015092510884566
$ws.Range("N23").Value = 1.007286299839288
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.005336274157109
$ws.Range("D24").Value = 1.026961740599293
$ws.Range("E24").Value = 1.007713048403275
$ws.Range("F24").Value = 1.017076468547901
$ws.Range("I24").Value = 1.029896203631762
$ws.Range("J24").Value = 1.011739425851277
$ws.Range("K24").Value = 1.030363648549301
$ws.Range("L24").Value = 1.011184883871151
$ws.Range("M24").Value = 1.020513922950725
$ws.Range("N24").Value = 1.01317621188412
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.012656933568266
$ws.Range("D25").Value = 1.028026619477224
$ws.Range("E25").Value = 1.014158177687067
$ws.Range("F25").Value = 1.023539203625232
$ws.Range("I25").Value = 1.03044354334937
$ws.Range("J25").Value = 1.018392760721128
$ws.Range("K25").Value = 1.031103780366417
$ws.Range("L25").Value = 1.01728031133148
$ws.Range("M25").Value = 1.026630770301134
$ws.Range("N25").Value = 1.019838995252635
